$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Use Case"

# Set explicit row heights (matches default 12.75) for rows 1 and 3
$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(3).RowHeight = 12.75

# Move the selection to A6
$ws.Range("A6").Select()
